$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank "spacer" cells in columns B, D and G (rows 1-7) that
# previously held only formatting with no content.
$ws.Range("B1:B7").Clear()
$ws.Range("D1:D7").Clear()
$ws.Range("G1:G7").Clear()

# Rows 6 and 7 lose their Index/Greeks figures (columns H:L).
$ws.Range("H6:L7").Clear()

# Replace the old Deal IDs with the new batch, and add a new deal on row 8.
$ws.Range("A2").Value = "SP2D01JL"
$ws.Range("A4").Value = "SP8T01PF"
$ws.Range("A3").Value = "SP8T01PG"
$ws.Range("A5").Value = "SP2D01JN"
$ws.Range("A6").Value = "SP1901P6"
$ws.Range("A7").Value = "SP8T01PJ"
$ws.Range("A8").Value = "SP8T01PK"

# Update the active selection to the newly added cell.
$null = $ws.Range("A8").Select()
